$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date in column C was updated from 2023-10-22
# (serial 45221) to 2023-10-25 (serial 45224) for all data rows (2-23).
$ws.Range("C2:C23").Value = 45224
